$p = $ppt.ActivePresentation

# --- Slide 12: "Component properties (Basic)" title run ---
# Only formatting (dirty="0") changes here; text itself is unchanged.
$s12 = $p.Slides.Item(12)
$titleShape = $s12.Shapes.Item(2)
$titleShape.TextFrame.TextRange.Runs(1).Font.Bold = $titleShape.TextFrame.TextRange.Runs(1).Font.Bold

# --- Slide 16: split "Everything requires a unique name" into two runs ---
$s16 = $p.Slides.Item(16)
$bodyShape = $s16.Shapes.Item(2)
$tr = $bodyShape.TextFrame.TextRange.Paragraphs(1, 1)
$tr.Text = "Everything requires a unique name"
$boldPart = $tr.Characters(23, 11)
$boldPart.Font.Bold = $true
